# The commit swaps the DrawingML colour scheme that lives in ppt/theme/theme1.xml
# (the "Integral" theme used by the one SlideMaster in this deck) for the stock
# "Office" colour scheme that used to live in ppt/theme/theme2.xml. The font
# scheme and format scheme (fills/lines/effects) are already byte-identical
# between the two theme parts, so only the 12 scheme colours actually change.
#
# PowerPoint's COM model exposes those 12 DrawingML colours via
# Slide.ThemeColorScheme.Colors(msoThemeColorSchemeIndex).RGB - note RGB is a
# Windows COLORREF (0x00BBGGRR), so values below are given as decimal BGR.
#
#   index  slot       new ("Office") hex   decimal BGR
#   1      dk1        000000               0
#   2      lt1        FFFFFF               16777215
#   3      dk2        44546A               6968388
#   4      lt2        E7E6E6               15132391
#   5      accent1    5B9BD5               13998939
#   6      accent2    ED7D31               3243501
#   7      accent3    A5A5A5               10855845
#   8      accent4    FFC000               49407
#   9      accent5    4472C4               12874308
#   10     accent6    70AD47               4697456
#   11     hlink      0563C1               12673797
#   12     folHlink   954F72               7491477

$p = $ppt.ActivePresentation

$officeThemeColorsBgr = @(0, 16777215, 6968388, 15132391, 13998939, 3243501, 10855845, 49407, 12874308, 4697456, 12673797, 7491477)

# All slides in this deck share the single SlideMaster/theme, so updating the
# scheme through the first slide updates ppt/theme/theme1.xml for everyone.
$slide = $p.Slides.Item(1)
$themeColors = $slide.ThemeColorScheme

for ($i = 1; $i -le $themeColors.Count; $i++) {
    $themeColors.Colors($i).RGB = $officeThemeColorsBgr[$i - 1]
}
